$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $r = $ws.Range($CellRef)
    if ($Value -match '^-?[0-9]*\.?[0-9]+([eE][-+]?[0-9]+)?$') {
        # Numeric-looking text: force literal text storage (quote-prefix),
        # matching the source workbook's inline-string-as-text cells, then
        # strip the auto-applied "quote prefix" style back to Normal so no
        # visible formatting changes.
        $r.Value = "'" + $Value
        $r.Style = "Normal"
    } else {
        $r.Value = $Value
    }
}

Set-CellText 'D2' '269.59'
Set-CellText 'D3' '22.67'
Set-CellText 'D4' '6.330'
Set-CellText 'D5' '0.06176'
Set-CellText 'D6' '3.649'
Set-CellText 'D7' '6.670'
Set-CellText 'D8' '1.370'
Set-CellText 'D9' '0.8290'
Set-CellText 'D10' '0.01373'
Set-CellText 'D11' '0.1607'
Set-CellText 'D12' '0.08283'
Set-CellText 'D13' '0.03557'
Set-CellText 'D14' '0.03244'
Set-CellText 'B15' 'ProBitToken'
Set-CellText 'C15' 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-CellText 'D15' '0.1239'
Set-CellText 'E15' '14ProBitTokenPROB'
Set-CellText 'B16' 'BitMartToken'
Set-CellText 'C16' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-CellText 'D16' '0.09329'
Set-CellText 'E16' '15BitMartTokenBMX'
Set-CellText 'B17' 'MCDex'
Set-CellText 'C17' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-CellText 'D17' '3.878'
Set-CellText 'E17' '16MCDexMCB'
Set-CellText 'B18' 'BitForexToken'
Set-CellText 'C18' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-CellText 'D18' '0.001653'
Set-CellText 'E18' '17BitForexTokenBF'
Set-CellText 'B19' 'CoinExToken'
Set-CellText 'C19' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-CellText 'D19' '0.04740'
Set-CellText 'E19' '18CoinExTokenCET'
Set-CellText 'B20' 'TigerCash'
Set-CellText 'C20' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-CellText 'D20' '0.006356'
Set-CellText 'E20' '19TigerCashTCH'
Set-CellText 'B21' 'HotbitToken'
Set-CellText 'C21' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-CellText 'D21' '0.005668'
Set-CellText 'E21' '20HotbitTokenHTB'
Set-CellText 'B22' 'BitKan'
Set-CellText 'C22' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-CellText 'D22' '0.001077'
Set-CellText 'E22' '21BitKanKAN'
Set-CellText 'B23' 'NitroEx'
Set-CellText 'C23' 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-CellText 'D23' '0.0001501'
Set-CellText 'E23' '22NitroExNTX'
Set-CellText 'B24' 'LEO'
Set-CellText 'C24' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-CellText 'D24' '3.731'
Set-CellText 'E24' '23LEOLEO'
Set-CellText 'B25' 'BTSEToken'
Set-CellText 'C25' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-CellText 'D25' '2.412'
Set-CellText 'E25' '24BTSETokenBTSE'
Set-CellText 'B26' 'BitpandaEcosystemToken'
Set-CellText 'C26' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-CellText 'D26' '0.3308'
Set-CellText 'E26' '25BitpandaEcosystemTokenBEST'
Set-CellText 'D27' '0.0002706'
Set-CellText 'D40' '0.04718'
Set-CellText 'D41' '0.006960'
Set-CellText 'D42' '0.1160'
Set-CellText 'D43' '0.003301'
Set-CellText 'E43' '42CEJICEJIWorstin24h'
Set-CellText 'D44' '0.01177'
Set-CellText 'D45' '0.00006259'
Set-CellText 'D46' '0.0009905'
Set-CellText 'D47' '0.00000000750'
Set-CellText 'D48' '0.9205'
Set-CellText 'E48' '47CoinbaseStockTokenCOIN'
Set-CellText 'D49' '0.002305'
Set-CellText 'D50' '0.00001401'
Set-CellText 'D51' '0.01241'
